$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '64.387.87'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.47%  '
$ws.Range("D3").Value = "'" + '3.401.27'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.26%  '
$ws.Range("D4").Value = "'" + '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = "'" + '580.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.33%  '
$ws.Range("D6").Value = "'" + '134.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.58%  '
$ws.Range("D7").Value = "'" + '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").Value = "'" + '3.402.78'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.04%  '
$ws.Range("D9").Value = "'" + '0.485'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.12%  '
$ws.Range("E10").Value = '  -10.91%  '
$ws.Range("E11").Value = '  -10.78%  '
$ws.Range("E12").Value = '  -7.89%  '
$ws.Range("D13").Value = "'" + '3.977.73'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.49%  '
$ws.Range("D14").Value = "'" + '0.0000175'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -10.47%  '
$ws.Range("E15").Value = '  -1.96%  '
$ws.Range("D16").Value = "'" + '25.99'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -8.31%  '
$ws.Range("D17").Value = "'" + '3.386.51'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.04%  '
$ws.Range("D18").Value = "'" + '64.405.83'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.61%  '
$ws.Range("D19").Value = "'" + '9.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -13.41%  '
$ws.Range("D20").Value = "'" + '5.76'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -6.46%  '
$ws.Range("E21").Value = '  -6.94%  '
$ws.Range("D22").Value = "'" + '377.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -9.45%  '
$ws.Range("D23").Value = "'" + '0.544'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -8.69%  '
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").Value = "'" + '71.84'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.73%  '
$ws.Range("D26").Value = "'" + '3.538.77'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.41%  '
$ws.Range("D27").Value = "'" + '0.0000103'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -9.11%  '
$ws.Range("D28").Value = "'" + '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").Value = "'" + '7.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -9.14%  '
$ws.Range("D30").Value = "'" + '2.17'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -11.41%  '
$ws.Range("D31").Value = "'" + '7.92'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -10.41%  '
$ws.Range("D32").Value = "'" + '3.416.68'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.19%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").Value = "'" + '22.79'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.64%  '
$ws.Range("D35").Value = "'" + '0.139'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -9.91%  '
$ws.Range("D36").Value = "'" + '167.85'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.37%  '
$ws.Range("D37").Value = "'" + '6.66'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -11.33%  '
$ws.Range("E38").Value = '  -11.71%  '
$ws.Range("D39").Value = "'" + '1.43'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.22%  '
$ws.Range("D40").Value = "'" + '4.57'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -12.32%  '
$ws.Range("D41").Value = "'" + '0.0744'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.94%  '
$ws.Range("D42").Value = "'" + '0.803'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.82%  '
$ws.Range("D43").Value = "'" + '1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("D44").Value = "'" + '41.99'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.96%  '
$ws.Range("D45").Value = "'" + '4.24'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -15.28%  '
$ws.Range("D46").Value = "'" + '1.58'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -10.42%  '
$ws.Range("E47").Value = '  +1.62%  '
$ws.Range("D48").Value = "'" + '22.39'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.38%  '
$ws.Range("D49").Value = "'" + '6.38'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.89%  '
$ws.Range("D50").Value = "'" + '2.148.68'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.58%  '
$ws.Range("D51").Value = "'" + '1.98'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -15.72%  '
